# Updated script on 21/11
$wb = $excel.ActiveWorkbook

# --- AdminSearchPage sheet updates ---
$ws = $wb.Worksheets.Item("AdminSearchPage")

# Fix the typo'd username value ("rtrtrtr" -> "nayana")
$ws.Range("A2").Value = "nayana"

# Usertype value changes from "staff" to "admin"
$ws.Range("B2").Value = "admin"

# Make AdminSearchPage the active/selected sheet and move the selection
# cursor to O16 (matches the saved view state of the workbook).
$ws.Activate() | Out-Null
$ws.Range("O16").Select() | Out-Null
